# Remove the trailing Jekyll-site boilerplate paragraphs that follow the
# "Requisitos" section: the blank paragraph, the "Ver no Jupiter Salvar em
# pdf Salvar em docx" line, and the "(c) 2020 ... Creative Commons
# Attribution" footer line. The paragraph holding the last real requirement
# text (LOM3073 ...) and the paragraphs that follow the removed block
# (blank paragraph + page-break paragraph) are left untouched.

$d = $word.ActiveDocument

$anchor = $null
$footer = $null

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*LOM3073*Requisito fraco*") {
        $anchor = $p
    }
    if ($p.Range.Text -like "*Contact: luizeleno@usp.br*") {
        $footer = $p
    }
}

if (($anchor -ne $null) -and ($footer -ne $null)) {
    $deleteRange = $d.Range($anchor.Range.End, $footer.Range.End)
    $deleteRange.Delete()
}
